$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C2:C180").Value = 45175
